# REPORTGEN-1102: part 1, added and removed counts missing when no previous snapshot selected
#
# For each "AxxChapter" sheet (A01..A10), the RepGen token in cell A3 that
# drives the "RULES_LIST_STATISTICS_RATIO" table needs an extra
# "EVOLUTION=true" flag appended, so the Added/Removed vulnerability counts
# are produced even when no previous snapshot is selected.
# The Summary sheet's "QUALITY_STANDARDS_EVOLUTION" token (cell B14) gets
# the same flag appended.

$wb = $excel.ActiveWorkbook

$chapterSheets = @("A01", "A02", "A03", "A04", "A05", "A06", "A07", "A08", "A09", "A10")

foreach ($sheetName in $chapterSheets) {
    $ws = $wb.Worksheets.Item($sheetName)
    $cell = $ws.Range("A3")
    $current = [string]$cell.Value2
    if ($current -like "RepGen:TABLE;RULES_LIST_STATISTICS_RATIO;*" -and $current -notlike "*EVOLUTION=true*") {
        $cell.Value2 = $current + ",EVOLUTION=true"
    }
}

$summary = $wb.Worksheets.Item("Summary")
$summaryCell = $summary.Range("B14")
$summaryCurrent = [string]$summaryCell.Value2
if ($summaryCurrent -like "RepGen:TABLE;QUALITY_STANDARDS_EVOLUTION;*" -and $summaryCurrent -notlike "*EVOLUTION=true*") {
    $summaryCell.Value2 = $summaryCurrent + ",EVOLUTION=true"
}

# Leave the A01 sheet's last selection on A3 (the cell that was edited).
$a01 = $wb.Worksheets.Item("A01")
[void]$a01.Activate()
[void]$a01.Range("A3").Select()

# Restore the originally selected / active cell on the Summary sheet.
[void]$summary.Activate()
[void]$summary.Range("B14").Select()
